$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(126).Insert()

$ws.Cells.Item(126, 1).Value = 11
$ws.Cells.Item(126, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(126, 3).Value = "Bíobío"
$ws.Cells.Item(126, 4).Value = 45119
$ws.Cells.Item(126, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(126, 5).Value = 8
$ws.Cells.Item(126, 6).Value = 100112024
$ws.Cells.Item(126, 7).Value = "Choclo"
$ws.Cells.Item(126, 8).Value = "Dulce o Americano"
$ws.Cells.Item(126, 9).Value = "Primera"
$ws.Cells.Item(126, 10).Value = 100
$ws.Cells.Item(126, 11).Value = 25000
$ws.Cells.Item(126, 12).Value = 26000
$ws.Cells.Item(126, 13).Value = 25500
$ws.Cells.Item(126, 14).Value = "$/malla 70 unidades"
$ws.Cells.Item(126, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(126, 16).Value = 364
$ws.Cells.Item(126, 17).Value = 70
$ws.Cells.Item(126, 18).Value = "Hortaliza"
